# Practiced Cambridge 19 Test1 -- fill in row 47 of Sheet1's score table.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C47").Value = 43
$ws.Range("D47").Value = 45531
$ws.Range("E47").Value = "IELTS19_Test1"
$ws.Range("F47").Value = 30
$ws.Range("G47").Formula = "=IFERROR(INDEX(Sheet2!`$F`$5:`$F`$20, MATCH(Table1[[#This Row],[Lis_Mark]], Sheet2!`$D`$5:`$D`$20, 1)),""No Grade"")"
$ws.Range("H47").Value = 35
$ws.Range("I47").Formula = "=IFERROR(INDEX(Sheet2!`$F`$5:`$F`$20, MATCH(Table1[[#This Row],[Read_Mark]], Sheet2!`$D`$5:`$D`$20, 1)),""No Grade"")"
$ws.Range("J47").Value = 6
$ws.Range("K47").Value = 4
$ws.Range("L47").Formula = "=(G47+I47+J47+K47)/4"

$ws.Range("N49").Select()

$wb.Save()
